# Prevent automatic recalculation so formula-cached values that are not
# part of this edit (e.g. the SUM() cells on the missing_values sheet)
# stay exactly as they were before - only the inputs actually touched by
# this edit should change.
$excel.Calculation = -4135  # xlCalculationManual

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "missing_values": update a handful of counts / percentages
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("missing_values")

$ws1.Range("B7").Value = 0
$ws1.Range("C7").Value = 0

$ws1.Range("B8").Value = 2161
$ws1.Range("C8").Value = 17.862456604397419

$ws1.Range("B18").Value = 0
$ws1.Range("C18").Value = 0

$ws1.Range("B19").Value = 165
$ws1.Range("C19").Value = 2.8032619775739041

# ---------------------------------------------------------------------
# Sheet "labor_incmon_imp_stochastic_reg": replace the summary row with
# the re-estimated statistics (and re-label it "2019" instead of the old
# "ila_monetario" tag).
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("labor_incmon_imp_stochastic_reg")

$a4 = $ws2.Range("A4")
$a4.NumberFormat = "@"
$a4.Value = "2019"
$a4.ClearFormats()

$ws2.Range("B4").Value = 2042580.7753590087
$ws2.Range("C4").Value = 251904.90625
$ws2.Range("D4").Value = 453428.84375
$ws2.Range("E4").Value = 1007619.5625
$ws2.Range("F4").Value = 2657232.5
$ws2.Range("G4").Value = 4042826.125
$ws2.Range("H4").Value = 2059723.0191452242
$ws2.Range("I4").Value = 251904.90625
$ws2.Range("J4").Value = 463618.8125
$ws2.Range("K4").Value = 1007619.625
$ws2.Range("L4").Value = 2653803.75
$ws2.Range("M4").Value = 4168662

# ---------------------------------------------------------------------
# New sheet "nonlabor_imp_stochastic_reg" with the non-labor imputation
# stats (two blocks of mean/p10/p25/p50/p75/p90/p99, without vs with
# imputation).
# ---------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws3 = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws3.Name = "nonlabor_imp_stochastic_reg"

$ws3.Range("B3").Value = "mean"
$ws3.Range("C3").Value = "p10"
$ws3.Range("D3").Value = "p25"
$ws3.Range("E3").Value = "p50"
$ws3.Range("F3").Value = "p75"
$ws3.Range("G3").Value = "p90"

$a4n = $ws3.Range("A4")
$a4n.NumberFormat = "@"
$a4n.Value = "2019"
$a4n.ClearFormats()

$ws3.Range("H3").Value = "p99"
$ws3.Range("I3").Value = "mean"
$ws3.Range("J3").Value = "p10"
$ws3.Range("K3").Value = "p25"
$ws3.Range("L3").Value = "p50"
$ws3.Range("M3").Value = "p75"
$ws3.Range("N3").Value = "p90"
$ws3.Range("O3").Value = "p99"

$ws3.Range("B4").Value = 434219.62010040088
$ws3.Range("C4").Value = 120914.3515625
$ws3.Range("D4").Value = 201523.921875
$ws3.Range("E4").Value = 300000
$ws3.Range("F4").Value = 478619.3125
$ws3.Range("G4").Value = 834161.9375
$ws3.Range("H4").Value = 2687251.75
$ws3.Range("I4").Value = 435681.67022345966
$ws3.Range("J4").Value = 115507.0625
$ws3.Range("K4").Value = 201523.921875
$ws3.Range("L4").Value = 300000
$ws3.Range("M4").Value = 490039.53125
$ws3.Range("N4").Value = 850000
$ws3.Range("O4").Value = 2683666.75

# Restore automatic calculation mode (doesn't retroactively recompute the
# stale SUM() caches above, it only affects calculation going forward).
$excel.Calculation = -4105  # xlCalculationAutomatic
